$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new header row at the top; existing rows (and the trailing
# blank row 78) shift down by one.
$ws.Rows.Item(1).Insert()

# The old trailing blank row (now row 79) is no longer needed once the
# header row takes row 1 - drop it so the sheet still ends at row 78.
$ws.Rows.Item(79).Delete()

# New header labels.
$ws.Range("A1").Value = "ISO Code"
$ws.Range("B1").Value = "Province"
$ws.Range("C1").Value = "Code"

# Clear the (now unused) explicit cell style from columns A and B so the
# whole data range falls back to the default style.
$ws.Range("A1:B78").ClearFormats()

# Column A gets an explicit width; column B keeps its existing width.
$ws.Columns.Item(1).ColumnWidth = 10.7

# View changes: zoom level and active selection.
$ws.Application.ActiveWindow.Zoom = 145
[void]$ws.Range("D3").Select()
